# Auto-generated edit script
# Updates numeric values in columns H-N across multiple leve-profit
# worksheets (ALC, ARM, BSM, CRP, CUL, GSM, WVR) to reflect refreshed
# market-board pricing data pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I40").Value = 2605.4546
$ws.Range("J40").Value = 2640
$ws.Range("K40").Value = 2605.4546
$ws.Range("L40").Value = 2640
$ws.Range("M40").Value = -2430.4546
$ws.Range("N40").Value = -2990
$ws.Range("H41").Value = 324.75
$ws.Range("I41").Value = 381.33334
$ws.Range("J41").Value = 300.5
$ws.Range("K41").Value = 381.33334
$ws.Range("L41").Value = 300.5
$ws.Range("M41").Value = 58.66665999999998
$ws.Range("N41").Value = -1180.5
$ws.Range("H64").Value = 3949.1091
$ws.Range("I64").Value = 3279.5454
$ws.Range("J64").Value = 6627.364
$ws.Range("K64").Value = 3279.5454
$ws.Range("L64").Value = 6627.364
$ws.Range("M64").Value = -3031.5454
$ws.Range("N64").Value = -7123.364
$ws.Range("H67").Value = 3949.1091
$ws.Range("I67").Value = 3279.5454
$ws.Range("J67").Value = 6627.364
$ws.Range("K67").Value = 3279.5454
$ws.Range("L67").Value = 6627.364
$ws.Range("M67").Value = -2421.5454
$ws.Range("N67").Value = -8343.364
$ws.Range("H76").Value = 3409.449
$ws.Range("I76").Value = 2981.5715
$ws.Range("J76").Value = 4479.143
$ws.Range("K76").Value = 2981.5715
$ws.Range("L76").Value = 4479.143
$ws.Range("M76").Value = -2666.5715
$ws.Range("N76").Value = -5109.143
$ws.Range("H79").Value = 3409.449
$ws.Range("I79").Value = 2981.5715
$ws.Range("J79").Value = 4479.143
$ws.Range("K79").Value = 2981.5715
$ws.Range("L79").Value = 4479.143
$ws.Range("M79").Value = -1889.5715
$ws.Range("N79").Value = -6663.143
$ws.Range("H87").Value = 15238.548
$ws.Range("J87").Value = 15238.548
$ws.Range("L87").Value = 15238.548
$ws.Range("N87").Value = -17734.548
$ws.Range("H90").Value = 15238.548
$ws.Range("J90").Value = 15238.548
$ws.Range("L90").Value = 45715.644
$ws.Range("N90").Value = -58195.644
$ws.Range("H116").Value = 2360
$ws.Range("I116").Value = 2068
$ws.Range("K116").Value = 2068
$ws.Range("M116").Value = 1374
$ws.Range("H137").Value = 1125.3334
$ws.Range("I137").Value = 1254.5454
$ws.Range("J137").Value = 922.2857
$ws.Range("K137").Value = 3763.6362
$ws.Range("L137").Value = 2766.8571
$ws.Range("M137").Value = -1213.6362
$ws.Range("N137").Value = -7866.8571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4243.78
$ws.Range("I32").Value = 4038.8572
$ws.Range("J32").Value = 6315.778
$ws.Range("K32").Value = 4038.8572
$ws.Range("L32").Value = 6315.778
$ws.Range("M32").Value = -3751.8572
$ws.Range("N32").Value = -6889.778
$ws.Range("H61").Value = 1951.5853
$ws.Range("I61").Value = 1923.4584
$ws.Range("J61").Value = 1991.2941
$ws.Range("K61").Value = 1923.4584
$ws.Range("L61").Value = 1991.2941
$ws.Range("M61").Value = -1711.4584
$ws.Range("N61").Value = -2415.2941
$ws.Range("H63").Value = 4408.3335
$ws.Range("I63").Value = 2612.5
$ws.Range("J63").Value = 8000
$ws.Range("K63").Value = 2612.5
$ws.Range("L63").Value = 8000
$ws.Range("M63").Value = -1926.5
$ws.Range("N63").Value = -9372
$ws.Range("H66").Value = 4408.3335
$ws.Range("I66").Value = 2612.5
$ws.Range("J66").Value = 8000
$ws.Range("K66").Value = 13062.5
$ws.Range("L66").Value = 40000
$ws.Range("M66").Value = -9630.5
$ws.Range("N66").Value = -46864
$ws.Range("H136").Value = 1951.5853
$ws.Range("I136").Value = 1923.4584
$ws.Range("J136").Value = 1991.2941
$ws.Range("K136").Value = 5770.3752
$ws.Range("L136").Value = 5973.8823
$ws.Range("M136").Value = -3220.3752
$ws.Range("N136").Value = -11073.8823

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -827
$ws.Range("N22").Value = -1346
$ws.Range("H35").Value = 34166.668
$ws.Range("I35").Value = 35000
$ws.Range("J35").Value = 34062.5
$ws.Range("K35").Value = 35000
$ws.Range("L35").Value = 34062.5
$ws.Range("M35").Value = -34690
$ws.Range("N35").Value = -34682.5
$ws.Range("H86").Value = 2658.25
$ws.Range("I86").Value = 2691.5833
$ws.Range("J86").Value = 2638.25
$ws.Range("K86").Value = 2691.5833
$ws.Range("L86").Value = 2638.25
$ws.Range("M86").Value = -1568.5833
$ws.Range("N86").Value = -4884.25
$ws.Range("H89").Value = 2658.25
$ws.Range("I89").Value = 2691.5833
$ws.Range("J89").Value = 2638.25
$ws.Range("K89").Value = 13457.9165
$ws.Range("L89").Value = 13191.25
$ws.Range("M89").Value = -7841.916499999999
$ws.Range("N89").Value = -24423.25
$ws.Range("H99").Value = 1538.4615
$ws.Range("I99").Value = 1000
$ws.Range("J99").Value = 1700
$ws.Range("K99").Value = 1000
$ws.Range("L99").Value = 1700
$ws.Range("M99").Value = 498
$ws.Range("N99").Value = -4696
$ws.Range("H105").Value = 2798.3076
$ws.Range("I105").Value = 2718.9092
$ws.Range("J105").Value = 3235
$ws.Range("K105").Value = 2718.9092
$ws.Range("L105").Value = 3235
$ws.Range("M105").Value = -971.9092000000001
$ws.Range("N105").Value = -6729

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4686.727
$ws.Range("I62").Value = 3500.5715
$ws.Range("J62").Value = 6762.5
$ws.Range("K62").Value = 3500.5715
$ws.Range("L62").Value = 6762.5
$ws.Range("M62").Value = -2876.5715
$ws.Range("N62").Value = -8010.5
$ws.Range("H65").Value = 4686.727
$ws.Range("I65").Value = 3500.5715
$ws.Range("J65").Value = 6762.5
$ws.Range("K65").Value = 17502.8575
$ws.Range("L65").Value = 33812.5
$ws.Range("M65").Value = -14382.8575
$ws.Range("N65").Value = -40052.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2339
$ws.Range("I109").Value = 508.9
$ws.Range("K109").Value = 1526.7
$ws.Range("M109").Value = -486.6999999999998
$ws.Range("H113").Value = 1076.05
$ws.Range("J113").Value = 829.5
$ws.Range("L113").Value = 2488.5
$ws.Range("N113").Value = -6828.5
$ws.Range("H122").Value = 924.1724
$ws.Range("I122").Value = 476.08
$ws.Range("J122").Value = 3724.75
$ws.Range("K122").Value = 4284.72
$ws.Range("L122").Value = 33522.75
$ws.Range("M122").Value = -1834.72
$ws.Range("N122").Value = -38422.75
$ws.Range("H131").Value = 850.0700000000001
$ws.Range("I131").Value = 200
$ws.Range("J131").Value = 870.1752300000001
$ws.Range("K131").Value = 600
$ws.Range("L131").Value = 2610.52569
$ws.Range("M131").Value = 4440
$ws.Range("N131").Value = -12690.52569

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5469.3657
$ws.Range("I70").Value = 4571.3335
$ws.Range("J70").Value = 6737.1763
$ws.Range("K70").Value = 4571.3335
$ws.Range("L70").Value = 6737.1763
$ws.Range("M70").Value = -4301.3335
$ws.Range("N70").Value = -7277.1763
$ws.Range("H73").Value = 5469.3657
$ws.Range("I73").Value = 4571.3335
$ws.Range("J73").Value = 6737.1763
$ws.Range("K73").Value = 4571.3335
$ws.Range("L73").Value = 6737.1763
$ws.Range("M73").Value = -3635.3335
$ws.Range("N73").Value = -8609.176299999999
$ws.Range("H80").Value = 3127.5
$ws.Range("I80").Value = 2920
$ws.Range("K80").Value = 2920
$ws.Range("M80").Value = -1922
$ws.Range("H83").Value = 3127.5
$ws.Range("I83").Value = 2920
$ws.Range("K83").Value = 14600
$ws.Range("M83").Value = -9608
$ws.Range("H113").Value = 1823.8
$ws.Range("I113").Value = 1635.3636
$ws.Range("K113").Value = 1635.3636
$ws.Range("M113").Value = 534.6364000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 90002000
$ws.Range("I13").Value = 135000500
$ws.Range("K13").Value = 135000500
$ws.Range("M13").Value = -135000360
